# Season-record columns: Wins / Losses / Ties, added after the existing
# "Unnamed: 28" column (AC) so the new columns land in AD:AF.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold font, border, centered/top alignment)
# from an existing header cell so the new header cells share the same
# style as the rest of row 1, then overwrite with the new header text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row shares the same team season record.
for ($r = 2; $r -le 56; $r++) {
    $ws.Cells.Item($r, 30).Value = 107
    $ws.Cells.Item($r, 31).Value = 55
    $ws.Cells.Item($r, 32).Value = 0
}

Write-Output "done"
